$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$words = @(
    'sprengen',
    'altern',
    'schreiten',
    'filmen',
    'formen',
    'schenken',
    'spüren',
    'tropfen',
    'treiben',
    'machen',
    'erben',
    'gründen',
    'drehen',
    'zeigen',
    'sperren',
    'sinken',
    'sorgen',
    'siegen',
    'tollen',
    'kehren',
    'warnen',
    'zielen',
    'lügen',
    'backen',
    'bitten',
    'betteln',
    'schmecken',
    'mögen',
    'kosten',
    'liefern',
    'fischen',
    'zünden',
    'spinnen',
    'feiern',
    'werfen',
    'lesen',
    'ändern',
    'brauchen',
    'malen',
    'kümmern',
    'hauen',
    'trauen',
    'schwingen',
    'wehtun',
    'fließen',
    'klingen',
    'heilen',
    'mauern',
    'loben',
    'spielen',
    'äußern',
    'sichern',
    'ärgern',
    'heulen',
    'folgen',
    'trennen',
    'jagen',
    'geben',
    'scheitern',
    'wundern',
    'klettern',
    'arten',
    'ehren',
    'planen',
    'decken',
    'saufen',
    'schlucken',
    'stecken',
    'bergen',
    'fallen',
    'graben',
    'knarren',
    'pflanzen',
    'dringen',
    'bauen',
    'schrecken',
    'führen',
    'schulden',
    'sterben',
    'stehlen',
    'zögern',
    'stammen',
    'flehen',
    'hören',
    'wüten',
    'schwächen',
    'seufzen',
    'schwören',
    'suchen',
    'scheinen',
    'gelten',
    'bellen',
    'räumen',
    'achten',
    'boxen',
    'wachsen',
    'biegen',
    'werden',
    'fangen',
    'lockern',
    'rufen',
    'wenden',
    'rasen',
    'reizen',
    'dienen',
    'liegen',
    'jubeln',
    'münzen',
    'töten',
    'greifen',
    'kichern',
    'enden',
    'quälen',
    'grüßen',
    'weichen',
    'pfeifen',
    'heben',
    'irren',
    'runden',
    'helfen',
    'flüchten',
    'freuen',
    'streichen',
    'bluten',
    'platzen',
    'wirken',
    'fällen',
    'fahren'
)

for ($i = 0; $i -lt $words.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 4).Value = $words[$i]
}
